$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Generate Report for Handback
#
# For the "f4248e7d-7a25-40ab-936f-957e9f12cc82" row (row 8) on both the
# zh-cn and de-de sheets, fill in the "Latest Target File", "Latest
# Handback File", "Latest Handback DateTime" and "Error Detail" columns,
# add a hyperlink on the new "Latest Target File" cell, and widen columns
# I (Latest Target File) and P (Error Detail) to 40 characters.
# -----------------------------------------------------------------

$targetFileName = "f4248e7d-7a25-40ab-936f-957e9f12cc82.md"
$targetFileUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a66f8f7a9ad107838d62d8bd5172e359c282a443/e2e/f4248e7d-7a25-40ab-936f-957e9f12cc82.md"
$errorDetail    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a894745aad697f4c0740b0032f8c4f8446968c62/e2e/f4248e7d-7a25-40ab-936f-957e9f12cc82.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a66f8f7a9ad107838d62d8bd5172e359c282a443/e2e/f4248e7d-7a25-40ab-936f-957e9f12cc82.md."

function Update-HandbackRow {
    param($ws, $handbackFile, $handbackDateTime)

    # Latest Target File (I8) + hyperlink
    $ws.Range("I8").Value = $targetFileName
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetFileUrl, "", "", $targetFileName)
    $ws.Range("I8").Font.Underline = $true
    $ws.Range("I8").Font.Color = 15570276

    # Latest Handback File (J8)
    $ws.Range("J8").Value = $handbackFile

    # Latest Handback DateTime (K8)
    $ws.Range("K8").Value = $handbackDateTime

    # Error Detail (P8)
    $ws.Range("P8").Value = $errorDetail

    # Widen columns I and P to 40 characters
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "f4248e7d-7a25-40ab-936f-957e9f12cc82.072ff206c711835003ab2d87a0dcccf7673d34b8.zh-cn.xlf" "2016-09-04 10:47:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "f4248e7d-7a25-40ab-936f-957e9f12cc82.072ff206c711835003ab2d87a0dcccf7673d34b8.de-de.xlf" "2016-09-04 10:47:38"
